$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Topics")
$ws.Activate()

# Insert 3 fresh rows right before the old "AOP" row (row 118), pushing the
# old rows 118-124 down to 121-127. Rows 113-117 do not exist yet in the
# original sheet, so we only need to make room for the extra 3 rows that
# are needed beyond the already-unused 113-117 gap.
$ws.Rows("118:120").Insert()

# Apply style 16 (wrap-text, used for the long multi-line link cells) to
# D117 up front, mirroring the look of D103.
$ws.Range("D103").Copy()
$ws.Range("D117").PasteSpecial(-4122)
$ws.Rows("117:117").RowHeight = 28.8

# Apply style 17 (section-header shading) to the merged A114:E114 banner,
# mirroring the look of A112:E112, then merge it.
$ws.Range("A112:E112").Copy()
$ws.Range("A114:E114").PasteSpecial(-4122)
$ws.Range("A114:E114").Merge()

# Fill in the new cell contents. The order below matches the order the
# values were actually typed in, which drives the shared-string table order.
$ws.Range("C117").Value = "Rest Assured"
$ws.Range("D118").Value = "Cucumber"
$ws.Range("D119").Value = "gherkin"
$ws.Range("D115").Value = "https://www.baeldung.com/junit-5-repeated-test"
$ws.Range("A114").Value = "Session 26"
$ws.Range("D116").Value = "https://www.baeldung.com/spring-boot-h2-database"
$ws.Range("C115").Value = "jUnit repeat test"
$ws.Range("C116").Value = "H2 Database in memory"
$ws.Range("D117").Value = "https://rest-assured.io/`nhttps://maven.apache.org/surefire/maven-surefire-plugin/"
$ws.Range("D113").Value = "Junit 5 and Mockito"

$ws.Range("A114:E114").Select()
